$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "38.433.73"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "2.078.68"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws.Range("D5") "228.70"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  +0.26%  "
Set-TextValue $ws.Range("D7") "60.21"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue $ws.Range("D9") "0.382"
$ws.Range("E10").Value = "  +0.57%  "
Set-TextValue $ws.Range("D11") "0.104"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "2.387.21"
$ws.Range("E12").Value = "  +2.14%  "
Set-TextValue $ws.Range("D13") "14.86"
$ws.Range("E13").Value = "  +1.84%  "
Set-TextValue $ws.Range("D14") "22.35"
$ws.Range("E14").Value = "  +5.62%  "
Set-TextValue $ws.Range("D15") "0.780"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "2.078.66"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "38.340.02"
$ws.Range("E18").Value = "  +1.74%  "
Set-TextValue $ws.Range("D19") "71.13"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("E21").Value = "  +1.05%  "
Set-TextValue $ws.Range("D22") "225.02"
$ws.Range("E23").Value = "  -0.17%  "
Set-TextValue $ws.Range("D24") "2.39"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("E25").Value = "  +2.77%  "
Set-TextValue $ws.Range("D26") "169.89"
$ws.Range("E26").Value = "  +1.26%  "
Set-TextValue $ws.Range("D27") "9.40"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +6.74%  "
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("E30").Value = "  +8.27%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  +5.01%  "
$ws.Range("E33").Value = "  +6.55%  "
Set-TextValue $ws.Range("D34") "4.49"
$ws.Range("E34").Value = "  +2.68%  "
Set-TextValue $ws.Range("D35") "0.0605"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("E39").Value = "  +0.14%  "
Set-TextValue $ws.Range("D40") "18.31"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "1.539.33"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("E45").Value = "  -0.97%  "
Set-TextValue $ws.Range("D46") "7.71"
$ws.Range("E46").Value = "  +9.14%  "
$ws.Range("E47").Value = "  +0.63%  "
Set-TextValue $ws.Range("D48") "4.08"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("D51").Value = "2.276.12"
$ws.Range("E51").Value = "  +2.15%  "
